$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 5 ("Workflow draft"): add a trailing empty paragraph after
# the "To-do: ... Feature store on AWS?" line.
# -----------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(2)
$tr5 = $shape5.TextFrame.TextRange
$null = $tr5.InsertAfter("`r")

# -----------------------------------------------------------------
# Slide 6 ("CI"): resize/reposition the content placeholder, turn on
# shrink-text-on-overflow, extend the Makefile paragraph, and append
# the "Implemented for" / "Questions" paragraphs.
# -----------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shape6 = $s6.Shapes.Item(2)

# Reposition / resize the shape (values below are EMU/12700 = points).
$shape6.Left = 66
$shape6.Top = 115
$shape6.Width = 828
$shape6.Height = 396.25

# Shrink text on overflow (writes <a:normAutofit/> into bodyPr).
$shape6.TextFrame.AutoSize = 2

$tr6 = $shape6.TextFrame.TextRange

# Merge the tail of paragraph 2 ("... decides the " + "environment
# variables which " + "is fed into the recipe") into a single run
# with extra wording, leaving the "Makefile" run untouched.
$full6 = $tr6.Text
$oldTail = " is the recipe, while the YAML file is the skeleton while pulls from the recipe + decides the environment variables which is fed into the recipe"
$newTail = " is the recipe, while the YAML file is the skeleton while pulls from the recipe + decides the environment variables which is fed into the recipe"
$idx = $full6.IndexOf(" is the recipe")
$tailRange = $tr6.Characters($idx + 1, $oldTail.Length)
$tailRange.Text = $newTail

# Append the new paragraphs at the end of the text body.
$null = $tr6.InsertAfter("`r`rImplemented for: Model prediction`rQuestions: How to implement it for modelling? To store data on AWS, then pull from there? Maybe leave this to the retraining phase")

# Bold "Implemented for" / "Questions" labels.
$full6b = $tr6.Text
$impIdx = $full6b.IndexOf("Implemented for")
$impRange = $tr6.Characters($impIdx + 1, "Implemented for".Length)
$impRange.Font.Bold = $true

$qIdx = $full6b.IndexOf("Questions")
$qRange = $tr6.Characters($qIdx + 1, "Questions".Length)
$qRange.Font.Bold = $true
